# Insert a new client row ("ALTAMIRANO VILLAVICENCIO JUAN ALEJANDRO") under
# asesor "LINDAO ZUÑIGA BRYAN JOSE" at row 170 in both the "VENTAS POR GRUPO"
# and "VENTA MENSUAL" sheets. This pushes all subsequent rows down by one,
# so the trailing totals row moves from 303 to 304 and its "X de 301" labels
# become "X de 302".

$wb = $excel.ActiveWorkbook

# ---- Sheet "VENTAS POR GRUPO" (columns A:R) ----
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Rows.Item(170).Insert()
$ws1.Cells.Item(170, 1).Value = "LINDAO ZUÑIGA BRYAN JOSE"
$ws1.Cells.Item(170, 2).Value = "ALTAMIRANO VILLAVICENCIO JUAN ALEJANDRO"
for ($c = 3; $c -le 18; $c++) {
    $ws1.Cells.Item(170, $c).Value = 0
}

# Update the "X de 301" -> "X de 302" summary labels on the (now) last row.
for ($c = 3; $c -le 18; $c++) {
    $cell = $ws1.Cells.Item(304, $c)
    $old = $cell.Value2
    $cell.Value = $old.Replace("de 301", "de 302")
}

# ---- Sheet "VENTA MENSUAL" (columns A:G) ----
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Rows.Item(170).Insert()
$ws2.Cells.Item(170, 1).Value = "LINDAO ZUÑIGA BRYAN JOSE"
$ws2.Cells.Item(170, 2).Value = "ALTAMIRANO VILLAVICENCIO JUAN ALEJANDRO"
for ($c = 3; $c -le 7; $c++) {
    $ws2.Cells.Item(170, $c).Value = 0
}
